$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGM")

# Row 13 (A13:F13 already carry style s="1" from the template; just fill values)
$ws.Range("A13").Value2 = "Records&Achievement"
$ws.Range("B13").Value2 = "AchievementScreen"
$ws.Range("C13").Value2 = "AchievementScreen"
$ws.Range("D13").Value2 = "MorseCode.wav"
$ws.Range("E13").Value2 = "S"
$ws.Range("F13").Value2 = "O"

# Row 14 (new row, no special style)
$ws.Range("A14").Value2 = "Records&Achievement"
$ws.Range("B14").Value2 = "ScoreScreen"
$ws.Range("C14").Value2 = "GameOver"
$ws.Range("D14").Value2 = "UfoLanding.wav"
$ws.Range("E14").Value2 = "S"
$ws.Range("F14").Value2 = "O"

# Row 15 (new row, no special style)
$ws.Range("A15").Value2 = "Records&Achievement"
$ws.Range("B15").Value2 = "HighScoreScreen"
$ws.Range("C15").Value2 = "HighScoreScreen"
$ws.Range("D15").Value2 = "UfoSounds.wav"
$ws.Range("E15").Value2 = "S"
$ws.Range("F15").Value2 = "O"

# Extend the FileName CONCAT formula down through the new rows as one shared formula block
$ws.Range("G13:G15").Formula = "=CONCAT(""BGM_"",A13,""_"",B13,""_"",C13,""_"",D13)"

# The BGM tab should become the selected/active tab (previously SFX was selected)
$ws.Activate()

Write-Output "edit complete"
